# Apply the cryptos-list price/volume refresh described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as plain text. Values that look like ordinary
# numbers (e.g. "14.00", "0.9984") would otherwise be auto-converted to
# numbers by Excel -- silently dropping trailing zeros / exact formatting --
# so those are entered with a leading apostrophe to force text, exactly as
# typing them in the Excel UI would. Values that already look non-numeric
# (e.g. "26.576.99", with two dots) are assigned as plain text.

$ws.Range("D2").Value = '26.576.99'
$ws.Range("E2").Value = '  +0.15%  '

$ws.Range("D3").Value = '1.853.15'
$ws.Range("E3").Value = '  +0.13%  '

$ws.Range("D4").Value = "'" + '0.9984'
$ws.Range("E4").Value = '  -0.68%  '

$ws.Range("D5").Value = "'" + '265.85'
$ws.Range("E5").Value = '  +2.74%  '

$ws.Range("D6").Value = "'" + '0.9985'
$ws.Range("E6").Value = '  -0.61%  '

$ws.Range("D7").Value = "'" + '0.5234'
$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("D8").Value = "'" + '0.3288'
$ws.Range("E8").Value = '  +0.26%  '

$ws.Range("D9").Value = "'" + '0.06819'
$ws.Range("E9").Value = '  +1.49%  '

$ws.Range("D10").Value = "'" + '18.89'
$ws.Range("E10").Value = '  -2.50%  '

$ws.Range("D11").Value = "'" + '0.7807'
$ws.Range("E11").Value = '  +1.17%  '

$ws.Range("D12").Value = "'" + '0.07790'
$ws.Range("E12").Value = '  +1.20%  '

$ws.Range("D13").Value = '1.865.46'
$ws.Range("E13").Value = '  +0.20%  '

$ws.Range("D14").Value = "'" + '88.55'
$ws.Range("E14").Value = '  -0.26%  '

$ws.Range("D15").Value = "'" + '5.030'
$ws.Range("E15").Value = '  -0.24%  '

$ws.Range("E16").Value = '  -0.70%  '

$ws.Range("D17").Value = "'" + '14.00'
$ws.Range("E17").Value = '  -1.32%  '

$ws.Range("D18").Value = "'" + '0.000007985'
$ws.Range("E18").Value = '  +1.35%  '

$ws.Range("D19").Value = "'" + '0.9997'
$ws.Range("E19").Value = '  -0.24%  '

$ws.Range("D20").Value = '26.598.37'
$ws.Range("E20").Value = '  +0.01%  '

$ws.Range("D21").Value = '2.083.17'
$ws.Range("E21").Value = '  +0.79%  '

$ws.Range("D22").Value = "'" + '4.657'
$ws.Range("E22").Value = '  +0.91%  '

$ws.Range("D23").Value = "'" + '9.571'
$ws.Range("E23").Value = '  -1.68%  '

$ws.Range("D24").Value = "'" + '6.005'
$ws.Range("E24").Value = '  +0.46%  '

$ws.Range("D25").Value = "'" + '144.58'
$ws.Range("E25").Value = '  +0.06%  '

$ws.Range("D26").Value = "'" + '2.243'
$ws.Range("E26").Value = '  -4.76%  '

$ws.Range("D27").Value = "'" + '1.662'
$ws.Range("E27").Value = '  +0.08%  '

$ws.Range("D28").Value = "'" + '17.07'
$ws.Range("E28").Value = '  +0.23%  '

$ws.Range("D29").Value = "'" + '112.45'
$ws.Range("E29").Value = '  +0.90%  '

$ws.Range("E30").Value = '  -0.13%  '

$ws.Range("E31").Value = '  -0.73%  '

$ws.Range("D32").Value = "'" + '0.08761'
$ws.Range("E32").Value = '  -0.25%  '

$ws.Range("D33").Value = "'" + '0.04848'
$ws.Range("E33").Value = '  -0.56%  '

$ws.Range("E34").Value = '  +0.29%  '

$ws.Range("D35").Value = "'" + '0.7211'
$ws.Range("E35").Value = '  +1.88%  '

$ws.Range("D36").Value = "'" + '2.849'
$ws.Range("E36").Value = '  -1.43%  '

$ws.Range("D37").Value = "'" + '3.106'
$ws.Range("E37").Value = '  -1.15%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = "'" + '2.241'
$ws.Range("E38").Value = '  +1.06%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = "'" + '0.01786'
$ws.Range("E39").Value = '  -1.48%  '

$ws.Range("D40").Value = "'" + '0.4908'
$ws.Range("E40").Value = '  -0.56%  '

$ws.Range("D41").Value = "'" + '0.9220'
$ws.Range("E41").Value = '  +1.68%  '

$ws.Range("D42").Value = "'" + '111.79'
$ws.Range("E42").Value = '  -2.45%  '

$ws.Range("D43").Value = "'" + '6.085'
$ws.Range("E43").Value = '  +0.04%  '

$ws.Range("D44").Value = "'" + '0.9978'
$ws.Range("E44").Value = '  -0.65%  '

$ws.Range("D45").Value = "'" + '7.771'
$ws.Range("E45").Value = '  +0.15%  '

$ws.Range("D46").Value = "'" + '0.4204'
$ws.Range("E46").Value = '  -2.05%  '

$ws.Range("E47").Value = '  +0.44%  '

$ws.Range("D48").Value = "'" + '9.137'
$ws.Range("E48").Value = '  -0.35%  '

$ws.Range("E49").Value = '  -3.38%  '

$ws.Range("D50").Value = "'" + '35.10'
$ws.Range("E50").Value = '  -0.92%  '

$ws.Range("D51").Value = "'" + '0.8934'
$ws.Range("E51").Value = '  +3.15%  '
